$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (open_price, close_price, high_price, low_price, shares_outstanding)
# fixed_ticker (col I) is reset to "VNET" for every data row (2..43).
$priceData = @{
    2 = @(18.1299991607666, 20.56999969482422, 20.94000053405762, 16.84000015258789, 263895614)
    3 = @(20.54999923706055, 20.26000022888184, 21.02000045776367, 15.02999973297119, 263895614)
    4 = @(18.34000015258789, 20.04000091552734, 20.47999954223633, 17.80999946594238, 263895614)
    5 = @(20.59000015258789, 19.19000053405762, 20.88999938964844, 18.22999954223633, 263895614)
    6 = @(19.80999946594238, 20.18000030517578, 21.1200008392334, 19.63999938964844, 263895614)
    7 = @(10.10000038146973, 9.579999923706056, 10.47999954223633, 8.880000114440918, 263895614)
    8 = @(7.909999847412109, 7.309999942779541, 8.079999923706055, 6.769999980926514, 263895614)
    9 = @(7.059999942779541, 7.170000076293945, 7.550000190734863, 6.960000038146973, 263895614)
    10 = @(5.480000019073486, 5.409999847412109, 5.710000038146973, 5.28000020980835, 263895614)
    11 = @(5.010000228881836, 4.460000038146973, 5.03000020980835, 4.300000190734863, 263895614)
    12 = @(5.71999979019165, 7.510000228881836, 7.940000057220459, 5.659999847412109, 263895614)
    13 = @(8.109999656677246, 9.109999656677246, 9.59000015258789, 7.96999979019165, 263895614)
    14 = @(6.920000076293945, 5.269999980926514, 7.179999828338623, 5.110000133514404, 263895614)
    15 = @(9.640000343322754, 9.4399995803833, 11.53999996185303, 8.399999618530273, 263895614)
    16 = @(10.11999988555908, 10.88000011444092, 11.47999954223633, 9.229999542236328, 263895614)
    17 = @(8.510000228881836, 8.979999542236328, 9.770000457763672, 8.239999771118164, 263895614)
    18 = @(8, 7.860000133514404, 8.229999542236328, 7, 263895614)
    19 = @(7.880000114440918, 7.360000133514404, 8.119999885559082, 7, 263895614)
    20 = @(7.590000152587891, 8.479999542236328, 8.800000190734863, 7.119999885559082, 263895614)
    21 = @(7.489999771118164, 9.239999771118164, 11.25, 7.260000228881836, 263895614)
    22 = @(13.35000038146973, 14.97000026702881, 17.3700008392334, 12.71000003814697, 263895614)
    23 = @(24, 21.73999977111816, 30.44000053405762, 21.04999923706055, 263895614)
    24 = @(23.51000022888184, 22.95999908447266, 26.57999992370605, 21.5, 263895614)
    25 = @(35, 37.72999954223633, 43.79999923706055, 32.09000015258789, 263895614)
    26 = @(33.52000045776367, 27.90999984741211, 34.88000106811523, 26.80999946594238, 263895614)
    27 = @(22.75, 17.32999992370605, 22.94000053405762, 14.10999965667725, 263895614)
    28 = @(17.18000030517578, 15.68000030517578, 19.65999984741211, 15.3100004196167, 263895614)
    29 = @(9.149999618530272, 9.920000076293944, 10.28999996185303, 7.730000019073486, 263895614)
    30 = @(6.440000057220459, 5.96999979019165, 7.28000020980835, 5.079999923706055, 263895614)
    31 = @(5.96999979019165, 5.159999847412109, 6.28000020980835, 4.610000133514404, 263895614)
    32 = @(5.550000190734863, 4.190000057220459, 6.019999980926514, 4.070000171661377, 263895614)
    33 = @(5.699999809265137, 5.880000114440918, 6.829999923706055, 5.699999809265137, 263895614)
    34 = @(3.740000009536743, 2.940000057220459, 4.190000057220459, 2.75, 263895614)
    35 = @(2.940000057220459, 2.970000028610229, 3.059999942779541, 2.539999961853028, 263895614)
    36 = @(3.099999904632568, 3.220000028610229, 3.359999895095825, 2.650000095367432, 263895614)
    37 = @(2.880000114440918, 1.700000047683716, 2.880000114440918, 1.629999995231628, 263895614)
    38 = @(1.549999952316284, 1.690000057220459, 1.899999976158142, 1.419999957084656, 263895614)
    39 = @(2.109999895095825, 2.069999933242798, 2.470000028610229, 1.840000033378601, 263895614)
    40 = @(4.139999866485596, 3.279999971389771, 4.5, 3.170000076293945, 263895614)
    41 = @(4.760000228881836, 6.96999979019165, 8.140000343322754, 4.610000133514404, 263895614)
    42 = @(8.069999694824219, 6.289999961853027, 8.229999542236328, 4.650000095367432, 263895614)
    43 = @(6.849999904632568, 8.3100004196167, 9.420000076293944, 6.409999847412109, 263895614)
}

foreach ($row in $priceData.Keys) {
    $vals = $priceData[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # open_price
    $ws.Cells.Item($row, 5).Value = $vals[1]   # close_price
    $ws.Cells.Item($row, 6).Value = $vals[2]   # high_price
    $ws.Cells.Item($row, 7).Value = $vals[3]   # low_price
    $ws.Cells.Item($row, 8).Value = $vals[4]   # shares_outstanding
    $ws.Cells.Item($row, 9).Value = "VNET"     # fixed_ticker
}
